$d = $word.ActiveDocument

$replacements = @(
    @("74×32=", "81×43="),
    @("75×19=", "61×73="),
    @("91×21=", "47×89="),
    @("91×96=", "27×64="),
    @("79×97=", "54×58="),
    @("52×24=", "94×90="),
    @("96×85=", "46×33="),
    @("58×96=", "94×44="),
    @("66×68=", "72×51="),
    @("33×99=", "23×13="),
    @("72×92=", "41×89="),
    @("48×92=", "90×67="),
    @("42×81=", "39×46="),
    @("79×90=", "73×92="),
    @("38×14=", "59×95="),
    @("56×39=", "70×46="),
    @("80×71=", "65×51="),
    @("75×57=", "29×41="),
    @("43×32=", "15×46="),
    @("87×15=", "40×51="),
    @("20×99=", "83×85="),
    @("62×44=", "83×49="),
    @("34×91=", "69×11="),
    @("45×11=", "38×63="),
    @("16×52=", "77×97=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
